$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.986.12'
$ws.Range('E2').Value = '  +5.87%  '
$ws.Range('D3').Value = '3.537.49'
$ws.Range('E3').Value = '  +8.65%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '190.10'
$ws.Range('E5').Value = '  +10.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '564.31'
$ws.Range('E6').Value = '  +6.31%  '
$ws.Range('D7').Value = '3.528.25'
$ws.Range('E7').Value = '  +8.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  +3.46%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.634'
$ws.Range('E10').Value = '  +4.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.151'
$ws.Range('E11').Value = '  +12.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.93'
$ws.Range('E12').Value = '  +2.98%  '
$ws.Range('E13').Value = '  +5.66%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.43'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = '4.106.93'
$ws.Range('E15').Value = '  +9.20%  '
$ws.Range('D16').Value = '3.543.41'
$ws.Range('E16').Value = '  +9.33%  '
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').Value = '67.056.33'
$ws.Range('E18').Value = '  +6.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.27'
$ws.Range('E19').Value = '  +5.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.05'
$ws.Range('E20').Value = '  +7.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.998'
$ws.Range('E21').Value = '  +3.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '432.98'
$ws.Range('E22').Value = '  +17.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.14'
$ws.Range('E23').Value = '  +9.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.32'
$ws.Range('E24').Value = '  +5.21%  '
$ws.Range('E25').Value = '  +2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.10'
$ws.Range('E26').Value = '  -1.41%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('E27').Value = '  +9.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.25'
$ws.Range('E28').Value = '  +8.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.16'
$ws.Range('E29').Value = '  +10.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.50'
$ws.Range('E30').Value = '  +6.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '644.61'
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.59'
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.75'
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '59.73'
$ws.Range('E35').Value = '  +4.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '38.57'
$ws.Range('E36').Value = '  +4.96%  '
$ws.Range('D37').Value = '0.0₃0813'
$ws.Range('E37').Value = '  +11.28%  '
$ws.Range('E38').Value = '  +17.97%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.391'
$ws.Range('E40').Value = '  +3.41%  '
$ws.Range('E41').Value = '  +14.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.16%  '
$ws.Range('D43').Value = '3.043.72'
$ws.Range('E43').Value = '  +4.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  +3.53%  '
$ws.Range('E45').Value = '  +10.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.36'
$ws.Range('E46').Value = '  +9.76%  '
$ws.Range('E47').Value = '  +5.72%  '
$ws.Range('E49').Value = '  +5.67%  '
$ws.Range('E50').Value = '  +5.45%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.65'
$ws.Range('E51').Value = '  +10.71%  '
